$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 288
$ws.Range("C3").Value = 178027
$ws.Range("C4").Value = 167983
$ws.Range("C8").Value = 64.81
